# Insert a new data row at row 96 (pushing the existing rows 96:213 down to
# 97:214) and populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("96:96").Insert()

$ws.Range("A96").Value2 = 4
$ws.Range("B96").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C96").Value2 = "Los Lagos"
$ws.Range("D96").Value2 = 44895
$ws.Range("E96").Value2 = 10
$ws.Range("F96").Value2 = "Fruta"
$ws.Range("G96").Value2 = 100103
$ws.Range("H96").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I96").Value2 = 100103004
$ws.Range("J96").Value2 = "Durazno"
$ws.Range("K96").Value2 = "Florida King"
$ws.Range("L96").Value2 = "Primera"
$ws.Range("M96").Value2 = 600
$ws.Range("N96").Value2 = 16000
$ws.Range("O96").Value2 = 17000
$ws.Range("P96").Value2 = 16500
$ws.Range("Q96").Value2 = "`$/bandeja 10 kilos empedrada"
$ws.Range("R96").Value2 = "Provincia de San Felipe de Aconcagua"
$ws.Range("S96").Value2 = 1650
$ws.Range("T96").Value2 = 10
